$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 3-6 (Lemari Baju, Hair Dryer, Vacum Culeaner, Gamis, One set dress... becomes just the single Mouse Robot row)
$ws.Rows("3:6").Delete()

# Update row 2 values to the new "Mouse Robot" product row
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "ELEK90"
$ws.Range("C2").Value = "Mouse Robot"
$ws.Range("D2").Value = 90000
$ws.Range("E2").Value = 130000

# Collapse the outline to level 1 for rows (was 5) while keeping columns at 4
$ws.Outline.ShowLevels(1, 4)

# Update the view: selection moves to E10, clearing the previous topLeftCell/selection
$ws.Range("E10").Select()

# Resize the saved window (best-effort; matches windowWidth change in workbook view)
$excel.ActiveWindow.Width = 19635
